# Hortaliza, Macroferia Regional de Talca - Coliflor
# A new daily price-report row is inserted before row 276, pushing every
# subsequent record down by one row (old row N -> new row N+1), and the
# final record of the sheet is duplicated into a new trailing row so the
# used range grows from A1:R392 to A1:R394.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new record: insert a blank row at 276, which shifts
#    the former rows 276:392 down to 277:393.
$ws.Rows("276:276").Insert()

# 2) Populate the freshly inserted row 276 with the new weekly record.
$ws.Range("A276").Value = 5
$ws.Range("B276").Value = "Macroferia Regional de Talca"
$ws.Range("C276").Value = "Maule"
$ws.Range("D276").Value = 45007
$ws.Range("E276").Value = 7
$ws.Range("F276").Value = 100112008
$ws.Range("G276").Value = "Coliflor"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 3000
$ws.Range("K276").Value = 900
$ws.Range("L276").Value = 1000
$ws.Range("M276").Value = 967
$ws.Range("N276").Value = "$/unidad"
$ws.Range("O276").Value = "Región del Maule"
$ws.Range("P276").Value = 967
$ws.Range("Q276").Value = 1
$ws.Range("R276").Value = "Hortaliza"

# 3) Append a duplicate of the (now shifted) last record as a new row 394.
$ws.Range("A393:R393").Copy($ws.Range("A394:R394"))
